$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

$ws.Range("B4").Value = "TryskoMys"
$ws.Range("B7").Value = "Naty338"

$ws.Range("C19").Select()
